$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("alpha3F")

# New row 16 mirrors row 15's formatting: copy A15 -> A16 (keeps the bold/
# bordered/centered style) and B15 -> B16 (keeps the "HexGrid-60degTilt5degRes"
# shared-string label, same as row 15), then overwrite the index value.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 14
$ws.Range("B15").Copy($ws.Range("B16"))

$ws.Range("C16").Value = 1.009158055596972
$ws.Range("D16").Value = 0.9467888361041986
$ws.Range("E16").Value = 1.006318362769948
$ws.Range("F16").Value = 1.009158055596972
$ws.Range("G16").Value = 0.9675593721181777
$ws.Range("H16").Value = 1.021198791407562
$ws.Range("I16").Value = 1.006318362769948
$ws.Range("J16").Value = 0.9467888361041986
$ws.Range("K16").Value = 0.9765535994370733
$ws.Range("L16").Value = 0.9928558275170226
$ws.Range("M16").Value = 0.9928902967944677
